# Lassa Genbank reference summary: "rerun lassa after adding back 1 paper"
# Adds one new reference row (row 4) to Sheet1, citing Clegg JC. 2002.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Add the new reference as row 4 --------------------------------------
# Fill order matches the order the strings were newly introduced into the
# shared-strings table (Authors, Journal, PMID, Viruses, Host, Country,
# GenBank, Gene, Comment, then Title last).

$ws.Cells.Item(4, 2).Value = "Clegg JC."                                                              # Authors
$ws.Cells.Item(4, 4).Value = "Curr Top Microbiol Immunol. 2002;262:1-24. doi: 10.1007/978-3-642-56029-3_1."  # Journal

# PMID needs to stay a text value (matches existing PMID column formatting/style)
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(4, 5).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "11987802"
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(4, 5).PasteSpecial(-4122)  # xlPasteFormats (restore shared style, keep text value)

$ws.Cells.Item(4, 14).Value = "LASV"                                     # Viruses
$ws.Cells.Item(4, 15).Value = "Homo Sapiens, rodent"                     # Host
$ws.Cells.Item(4, 17).Value = "Nigeria, Sierra Leone, and Liberia"       # Country
$ws.Cells.Item(4, 18).Value = "X52400"                                   # GenBank
$ws.Cells.Item(4, 22).Value = "G, N"                                     # Gene
$ws.Cells.Item(4, 23).Value = "Book"                                     # Comment

$ws.Cells.Item(4, 1).Value = "Molecular phylogeny of the arenaviruses"   # Title

$ws.Cells.Item(4, 3).Value = 2002                                        # Year
$ws.Cells.Item(4, 13).Value = 7                                          # NumSeqs
$ws.Cells.Item(4, 19).Value = "Sanger"                                   # SeqMethod
$ws.Cells.Item(4, 20).Value = "Yes"                                      # CloneMethod

$ws.Rows.Item(4).RowHeight = 32

$excel.CutCopyMode = $false

# --- View/selection bookkeeping -------------------------------------------
$excel.ActiveWindow.ScrollColumn = 14
[void]$ws.Range("Q10").Select()

$wb.Windows.Item(1).Left = 420
